$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").Value = "2023-12-06 11:02:32"
$ws.Range("B19").Value = 0.001

$ws.Range("A20").Value = "2023-12-06 11:02:51"
$ws.Range("B20").Value = 0.0008
